$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 287 (pushes the old rows 287..364 down to 288..365)
$ws.Rows(287).Insert()

# Populate the newly inserted row 287 with the new weekly price entry
$ws.Range("A287").Value = 11
$ws.Range("B287").Value = "Vega Monumental Concepción"
$ws.Range("C287").Value = "Bíobío"
$ws.Range("D287").Value = 44932
$ws.Range("E287").Value = 8
$ws.Range("F287").Value = 100112008
$ws.Range("G287").Value = "Coliflor"
$ws.Range("H287").Value = "Sin especificar"
$ws.Range("I287").Value = "Primera"
$ws.Range("J287").Value = 2000
$ws.Range("K287").Value = 750
$ws.Range("L287").Value = 800
$ws.Range("M287").Value = 775
$ws.Range("N287").Value = "`$/unidad"
$ws.Range("O287").Value = "Región Metropolitana"
$ws.Range("P287").Value = 775
$ws.Range("Q287").Value = 1
$ws.Range("R287").Value = "Hortaliza"
